$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q2" (i.e.
#    right before "总计"), and fill it with the quarter's holdings.
# ------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q2")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the header-row formatting (bold font + border) from the
# "2021-Q2" sheet so the new sheet's header matches the existing style.
$afterSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0

# Columns B..G hold numeric-looking text ("003366", "0.09", ...):
# force text entry so leading zeros / original formatting survive,
# then drop back to the default "Normal" style so no stray number
# format sticks around on the cell.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "003366"
$newSheet.Range("C2").Value = "浙商汇金中证转型成长指数"
$newSheet.Range("D2").Value = "0.09"
$newSheet.Range("E2").Value = "93.88"
$newSheet.Range("F2").Value = "1.21"
$newSheet.Range("G2").Value = "0.0011"
$newSheet.Range("B2:G2").Style = "Normal"

$newSheet.Range("H2").Value = 7

# Copy the "A" column formatting used for the row index cells.
$afterSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to "总计", above the existing rows.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2020-Q4"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.02

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.15

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0

# The former row 3 ("A" index cell) already carries the bold/border
# style; copy it down onto the newly written row 4's "A" cell (it
# previously had no style since row 4 didn't exist before).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$totalSheet.Range("A4").Value = 2
